# Auto-generated COM-interop script: insert 4 new slides
# ('Part V.1 ...' .. 'Part V.2 - Create a Zotero Snapshot')
# before the existing final 'Notes' slide (position 40-43,
# 'Notes' shifts from 40 -> 44), matching layout 'Title and
# Content' (slideLayout2.xml) used by the rest of the deck.

$p = $ppt.ActivePresentation

# --- New slide at position 40: Part V.1 - BatLit Release Introduction. ---
$s1 = $p.Slides.Add(40 + 0, 2)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = 'Part V.1 - BatLit Release Introduction.'

$s1Body = $s1.Shapes.Item(2).TextFrame.TextRange
$s1Body.Text = 'Metadata for BatLit releases are in the data folder of https://github.com/bat-literature/bat-literature.github.io . A new release needs to be linked to the older releases to enable version tracking between versions. This is why we need to clone or update the github repository first before making a snapshot of the BatLit Zotero group. Then, after making the snapshot, we include the tracked metadata and commit this to the repository. Also, we keep a copy of the metadata + pdfs elsewhere. So, you need to backup the data folder including metadata and pdfs.'


# --- New slide at position 41: Part V.2 - Clone/Update BatLit Repository ---
$s2 = $p.Slides.Add(40 + 1, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = 'Part V.2 - Clone/Update BatLit Repository'

$s2Body = $s2.Shapes.Item(2).TextFrame.TextRange
$s2Body.Text = 'Make sure to install git on your system by running:sudo apt install gitcheck whether git installed bygit --versionthis should produce something like:git version 2.43.0then, rungit clone https://github.com/bat-literature/bat-literature.github.io to “clone” (or create a copy of) the BatLit repository. By default, the repository is cloned into a folder with the same name as the repository (e.g., bat-literature.github.io).'

$s2Body.Paragraphs(1).Characters(22, 3).Font.Name = 'Courier'
$s2Body.Paragraphs(2).Font.Name = 'Courier'
$s2Body.Paragraphs(4).Font.Name = 'Courier'
$s2Body.Paragraphs(6).Font.Name = 'Courier'
$s2Body.Paragraphs(8).Font.Name = 'Courier'
$s2Body.Paragraphs(9).Characters(152, 24).Font.Name = 'Courier'

# --- New slide at position 42: Part V.2 - Verify Current BatLit Version ---
$s3 = $p.Slides.Add(40 + 2, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = 'Part V.2 - Verify Current BatLit Version'

$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange
$s3Body.Text = 'Step 1. go into bat-literature.github.io directory using cd Step 2. run preston head --algo md5 Step 3. Compare the output of the command with the existing version on https://batlit.org/datapaper'

$s3Body.Paragraphs(1).Characters(17, 24).Font.Name = 'Courier'
$s3Body.Paragraphs(1).Characters(58, 2).Font.Name = 'Courier'
$s3Body.Paragraphs(1).Characters(73, 23).Font.Name = 'Courier'

# --- New slide at position 43: Part V.2 - Create a Zotero Snapshot ---
$s4 = $p.Slides.Add(40 + 3, 2)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = 'Part V.2 - Create a Zotero Snapshot'

$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
$s4Body.Text = 'Step 1. Set Zotero API Key with read-only access to the BatLit Zotero Step 2. Go into the bat-literature.github.io directory Step 3. Verify that a preston history exists by running preston history --algo md5 Step 4. Create a snapshot of the BatLit Zotero group'

$s4Body.Paragraphs(1).Characters(91, 24).Font.Name = 'Courier'
$s4Body.Paragraphs(1).Characters(182, 26).Font.Name = 'Courier'

Write-Host "Inserted 4 slides; total slide count:" $p.Slides.Count
